$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(92, 4).Value = 44935
$ws.Cells.Item(92, 8).Value = "Americana (o)"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 150
$ws.Cells.Item(92, 11).Value = 12000
$ws.Cells.Item(92, 12).Value = 12000
$ws.Cells.Item(92, 13).Value = 12000
$ws.Cells.Item(92, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(92, 15).Value = "Región del Maule"
$ws.Cells.Item(92, 16).Value = 800
$ws.Cells.Item(92, 17).Value = 15

$ws.Cells.Item(93, 4).Value = 44419
$ws.Cells.Item(93, 8).Value = "Americana (o)"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 50
$ws.Cells.Item(93, 11).Value = 80000
$ws.Cells.Item(93, 12).Value = 80000
$ws.Cells.Item(93, 13).Value = 80000
$ws.Cells.Item(93, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(93, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(93, 16).Value = 3200
$ws.Cells.Item(93, 17).Value = 25

$ws.Cells.Item(94, 4).Value = 44419
$ws.Cells.Item(94, 8).Value = "Americana (o)"
$ws.Cells.Item(94, 9).Value = "Segunda"
$ws.Cells.Item(94, 10).Value = 50
$ws.Cells.Item(94, 11).Value = 70000
$ws.Cells.Item(94, 12).Value = 70000
$ws.Cells.Item(94, 13).Value = 70000
$ws.Cells.Item(94, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(94, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(94, 16).Value = 2800
$ws.Cells.Item(94, 17).Value = 25

$ws.Cells.Item(95, 4).Value = 44222
$ws.Cells.Item(95, 8).Value = "Americana (o)"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 100
$ws.Cells.Item(95, 11).Value = 13000
$ws.Cells.Item(95, 12).Value = 13000
$ws.Cells.Item(95, 13).Value = 13000
$ws.Cells.Item(95, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(95, 15).Value = "Región del Maule"
$ws.Cells.Item(95, 16).Value = 929
$ws.Cells.Item(95, 17).Value = 14

$ws.Cells.Item(96, 4).Value = 44222
$ws.Cells.Item(96, 8).Value = "Americana (o)"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 50
$ws.Cells.Item(96, 11).Value = 20000
$ws.Cells.Item(96, 12).Value = 20000
$ws.Cells.Item(96, 13).Value = 20000
$ws.Cells.Item(96, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(96, 15).Value = "Región del Maule"
$ws.Cells.Item(96, 16).Value = 800
$ws.Cells.Item(96, 17).Value = 25

$ws.Cells.Item(97, 4).Value = 44550
$ws.Cells.Item(97, 8).Value = "Americana (o)"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 100
$ws.Cells.Item(97, 11).Value = 17000
$ws.Cells.Item(97, 12).Value = 17000
$ws.Cells.Item(97, 13).Value = 17000
$ws.Cells.Item(97, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(97, 15).Value = "Región del Maule"
$ws.Cells.Item(97, 16).Value = 1214
$ws.Cells.Item(97, 17).Value = 14

$ws.Cells.Item(98, 4).Value = 44446
$ws.Cells.Item(98, 8).Value = "Americana (o)"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 50
$ws.Cells.Item(98, 11).Value = 75000
$ws.Cells.Item(98, 12).Value = 75000
$ws.Cells.Item(98, 13).Value = 75000
$ws.Cells.Item(98, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(98, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(98, 16).Value = 3000
$ws.Cells.Item(98, 17).Value = 25

$ws.Cells.Item(99, 4).Value = 44193
$ws.Cells.Item(99, 8).Value = "Americana (o)"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 100
$ws.Cells.Item(99, 11).Value = 15000
$ws.Cells.Item(99, 12).Value = 15000
$ws.Cells.Item(99, 13).Value = 15000
$ws.Cells.Item(99, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(99, 15).Value = "Región del Maule"
$ws.Cells.Item(99, 16).Value = 1071
$ws.Cells.Item(99, 17).Value = 14

$ws.Cells.Item(100, 4).Value = 44574
$ws.Cells.Item(100, 8).Value = "Americana (o)"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 150
$ws.Cells.Item(100, 11).Value = 15000
$ws.Cells.Item(100, 12).Value = 15000
$ws.Cells.Item(100, 13).Value = 15000
$ws.Cells.Item(100, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(100, 15).Value = "Región del Maule"
$ws.Cells.Item(100, 16).Value = 1071
$ws.Cells.Item(100, 17).Value = 14

$ws.Cells.Item(101, 4).Value = 44236
$ws.Cells.Item(101, 8).Value = "Americana (o)"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 150
$ws.Cells.Item(101, 11).Value = 12000
$ws.Cells.Item(101, 12).Value = 12000
$ws.Cells.Item(101, 13).Value = 12000
$ws.Cells.Item(101, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(101, 15).Value = "Región del Maule"
$ws.Cells.Item(101, 16).Value = 480
$ws.Cells.Item(101, 17).Value = 25

$ws.Cells.Item(102, 4).Value = 44315
$ws.Cells.Item(102, 8).Value = "Cacho cabra rojo"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 150
$ws.Cells.Item(102, 11).Value = 15000
$ws.Cells.Item(102, 12).Value = 15000
$ws.Cells.Item(102, 13).Value = 15000
$ws.Cells.Item(102, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(102, 15).Value = "Región del Maule"
$ws.Cells.Item(102, 16).Value = 600
$ws.Cells.Item(102, 17).Value = 25

$ws.Cells.Item(103, 4).Value = 44551
$ws.Cells.Item(103, 8).Value = "Americana (o)"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 200
$ws.Cells.Item(103, 11).Value = 15000
$ws.Cells.Item(103, 12).Value = 15000
$ws.Cells.Item(103, 13).Value = 15000
$ws.Cells.Item(103, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(103, 15).Value = "Región del Maule"
$ws.Cells.Item(103, 16).Value = 1071
$ws.Cells.Item(103, 17).Value = 14

$ws.Cells.Item(104, 4).Value = 44258
$ws.Cells.Item(104, 8).Value = "Cacho cabra verde"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 150
$ws.Cells.Item(104, 11).Value = 11000
$ws.Cells.Item(104, 12).Value = 11000
$ws.Cells.Item(104, 13).Value = 11000
$ws.Cells.Item(104, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(104, 15).Value = "Región del Maule"
$ws.Cells.Item(104, 16).Value = 440
$ws.Cells.Item(104, 17).Value = 25

$ws.Cells.Item(105, 4).Value = 44258
$ws.Cells.Item(105, 8).Value = "Cristal"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 150
$ws.Cells.Item(105, 11).Value = 12000
$ws.Cells.Item(105, 12).Value = 12000
$ws.Cells.Item(105, 13).Value = 12000
$ws.Cells.Item(105, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(105, 15).Value = "Región del Maule"
$ws.Cells.Item(105, 16).Value = 480
$ws.Cells.Item(105, 17).Value = 25

$ws.Cells.Item(106, 4).Value = 44266
$ws.Cells.Item(106, 8).Value = "Cacho cabra verde"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 150
$ws.Cells.Item(106, 11).Value = 12000
$ws.Cells.Item(106, 12).Value = 12000
$ws.Cells.Item(106, 13).Value = 12000
$ws.Cells.Item(106, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(106, 15).Value = "Región del Maule"
$ws.Cells.Item(106, 16).Value = 480
$ws.Cells.Item(106, 17).Value = 25

$ws.Cells.Item(107, 4).Value = 44266
$ws.Cells.Item(107, 8).Value = "Cristal"
$ws.Cells.Item(107, 9).Value = "Primera"
$ws.Cells.Item(107, 10).Value = 100
$ws.Cells.Item(107, 11).Value = 12000
$ws.Cells.Item(107, 12).Value = 12000
$ws.Cells.Item(107, 13).Value = 12000
$ws.Cells.Item(107, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(107, 15).Value = "Región del Maule"
$ws.Cells.Item(107, 16).Value = 480
$ws.Cells.Item(107, 17).Value = 25

$ws.Cells.Item(108, 4).Value = 44312
$ws.Cells.Item(108, 8).Value = "Cristal"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 200
$ws.Cells.Item(108, 11).Value = 14000
$ws.Cells.Item(108, 12).Value = 14000
$ws.Cells.Item(108, 13).Value = 14000
$ws.Cells.Item(108, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(108, 15).Value = "Región del Maule"
$ws.Cells.Item(108, 16).Value = 560
$ws.Cells.Item(108, 17).Value = 25

$ws.Cells.Item(109, 4).Value = 44298
$ws.Cells.Item(109, 8).Value = "Cacho cabra rojo"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 150
$ws.Cells.Item(109, 11).Value = 13000
$ws.Cells.Item(109, 12).Value = 13000
$ws.Cells.Item(109, 13).Value = 13000
$ws.Cells.Item(109, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(109, 15).Value = "Región del Maule"
$ws.Cells.Item(109, 16).Value = 520
$ws.Cells.Item(109, 17).Value = 25
